# TC47_Canine_Filter_Breed-YorkshireTerr.xlsx - "Fixed ICDC breed all testcases"
#
# The StatQuery (column D) Cypher query text used by CasesTab / SamplesTab /
# FilesTab all shares a single string, so updating any one of the three
# cells updates them all.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lines = @(
    "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)",
    "OPTIONAL MATCH (samp:sample)-->(c)",
    "OPTIONAL MATCH (diag:diagnosis)-->(c)",
    "OPTIONAL MATCH (f:file)-[*]->(c)",
    "OPTIONAL MATCH (sf:file)-->(s)",
    "WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p",
    "WHERE demo.breed IN ['Yorkshire Terrier']",
    "RETURN  ",
    "    count(distinct p) AS Programs,",
    "    count(distinct s) AS Studies,",
    "    count(distinct c) AS Cases,",
    "    count(distinct samp) AS Samples,",
    "    count(distinct f) AS " + [char]96 + "Case Files" + [char]96 + ",",
    "    count(distinct sf) AS " + [char]96 + "Study Files" + [char]96
)
$newQuery = [string]::Join([char]10, $lines)

$ws.Range("D2").Value = $newQuery
$ws.Range("D3").Value = $newQuery
$ws.Range("D4").Value = $newQuery

# The sheet's saved view scrolled down two rows (topLeftCell C1 -> C3) and
# the selection moved from C2 to C4.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 3
$ws.Range("C4").Select()
